# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Periodo Mora" detail table (rows 16-29) is rebuilt: instead of
# alternating between the two workers period-by-period, the rows are
# regrouped so that each worker's 7 periods are listed together, newest
# period (2206) first down to the oldest (2112).
#   rows 16-22 -> EMILSA ISABEL LOPEZ DE ORDOÑEZ (CC 26024057), periods 2206..2112
#   rows 23-29 -> CESAR JOSE ORDOÑEZ LOPEZ       (CC 73169249), periods 2206..2112
# The "Valor Mora" dates that used to sit on period 2206 (rows 28/29) now
# sit on the new first row of each worker's block (rows 16/23), and the
# other period rows keep the common date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cc1 = "26024057"
$name1 = "EMILSA ISABEL LOPEZ DE ORDOÑEZ"
$cc2 = "73169249"
$name2 = "CESAR JOSE ORDOÑEZ LOPEZ"

$periods = @("2206", "2205", "2204", "2203", "2202", "2201", "2112")

# Rows 16-22: worker 1 (EMILSA), newest period first
for ($i = 0; $i -lt 7; $i++) {
    $r = 16 + $i
    $ws.Cells.Item($r, 3).Value = $cc1
    $ws.Cells.Item($r, 4).Value = $name1
    $ws.Cells.Item($r, 5).Value = $periods[$i]
}
$ws.Cells.Item(16, 6).Value = 26650
for ($i = 1; $i -lt 7; $i++) {
    $ws.Cells.Item(16 + $i, 6).Value = 36341
}

# Rows 23-29: worker 2 (CESAR), newest period first
for ($i = 0; $i -lt 7; $i++) {
    $r = 23 + $i
    $ws.Cells.Item($r, 3).Value = $cc2
    $ws.Cells.Item($r, 4).Value = $name2
    $ws.Cells.Item($r, 5).Value = $periods[$i]
}
$ws.Cells.Item(23, 6).Value = 25749
for ($i = 1; $i -lt 7; $i++) {
    $ws.Cells.Item(23 + $i, 6).Value = 36341
}
